# Applies the cryptos.xlsx price/volume/coin-order refresh described by the diff.
# Each target cell stores plain text (coin names, URLs, price/volume strings), so
# values are set with a leading apostrophe to keep numeric-looking strings
# (e.g. "7.00", "0.0000133") stored as text instead of being auto-converted to
# numbers by Excel, then the cell Style is reset to "Normal" so the quote-prefix
# indicator does not leave a stray style index behind (matches the original,
# un-styled inline-string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.360.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.84%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.587.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.76%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'518.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'141.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.93%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.26%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.605.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.55%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.36%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.101"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.76%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.44%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.038.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.98%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'58.275.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.98%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'20.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'ShibaInu"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'WrappedEther"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'2.583.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.15%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'337.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.88%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.63%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.30%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.20%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.21%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'Polygon"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'0.401"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Binance-PegBSC-USD"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.03%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.683.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.55%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0734"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.57%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.65%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'EthereumClassic"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'18.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.19%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Monero"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'149.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.74%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.77%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Fetch.AI"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.01%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Stacks"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.14%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'36.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.13%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.828"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Filecoin"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'3.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.01%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'FirstDigitalUSD"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Bittensor"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'274.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.40%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Mantle"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.602"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'WhiteBITCoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'10.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.37%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0946"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.99%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'18.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Hedera"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.51%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.970.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.41%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'RenderToken"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'4.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.95%  "
$ws.Range("E51").Style = "Normal"
